$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data table (row 2) for the new year 2023.
# This shifts all existing data rows (2-24) down to (3-25) and the trailing
# blank filler rows (31-34) down to (32-35), copying formatting from the
# row above as Excel normally does on a row insert.
$ws.Rows(2).Insert()

# Populate the new 2023 row. Only Jan/Feb/Mar are known so far.
$ws.Range("A2").Value = 2023
$ws.Range("B2").Value = 0.4543
$ws.Range("C2").Value = 0.4552
$ws.Range("D2").Value = 0.3298

# The row-insert copied formatting (and therefore empty styled cells) into
# E2:M2 from the row above; remove them so the row only carries the cells
# that actually have data (plus the untouched, still-empty O2 carried over).
$ws.Range("E2:M2").Clear()

# The year 2022 row (shifted from row 2 to row 3) is now complete: fill in
# the December value that was previously missing.
$ws.Range("M3").Value = 0.3976
